$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$resetRow = @(27, 49, 61)

for ($row = 2; $row -le 95; $row++) {
    if ($row -eq 36) { continue }

    $eCell = $ws.Cells.Item($row, 5)
    $eVal = $eCell.Value2

    if ($resetRow -contains $row) {
        $dVal = $ws.Cells.Item($row, 4).Value2
        $fCell = $ws.Cells.Item($row, 6)
        $fVal = $fCell.Value2
        $eCell.Value2 = $dVal

        $dateStr = [Math]::Round($fVal).ToString()
        $dt = [DateTime]::ParseExact($dateStr, "yyyyMMdd", $null)
        $dt2 = $dt.AddDays($dVal)
        $fCell.Value2 = [int]$dt2.ToString("yyyyMMdd")
    } else {
        $eCell.Value2 = $eVal - 1
    }
}
